$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.8
$ws.Range("I3").Value = 2.8
$ws.Range("J3").Value = 3.6
$ws.Range("AJ3").Value = 12
$ws.Range("AQ3").Value = 51
$ws.Range("O5").Value = 1.62
$ws.Range("P5").Value = 2.2
$ws.Range("G6").Value = 1.53
$ws.Range("H6").Value = 3.75
$ws.Range("I6").Value = 7.5
$ws.Range("J6").Value = 2.2
$ws.Range("L6").Value = 7.5
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7
$ws.Range("Q6").Value = 2.4
$ws.Range("R6").Value = 1.53
$ws.Range("X6").Value = 6
$ws.Range("Z6").Value = 10
$ws.Range("AA6").Value = 17
$ws.Range("AI6").Value = 34
$ws.Range("AJ6").Value = 23
$ws.Range("AR6").Value = 67
$ws.Range("AS6").Value = 301
$ws.Range("AX6").Value = 8
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.48
$ws.Range("L8").Value = 4
$ws.Range("O8").Value = 1.62
$ws.Range("P8").Value = 2.2
$ws.Range("S8").Value = 1.67
$ws.Range("T8").Value = 2.1
$ws.Range("AT8").Value = 2.1
$ws.Range("BB8").Value = 126
$ws.Range("G13").Value = 7
$ws.Range("H13").Value = 4.33
$ws.Range("AG13").Value = 201
$ws.Range("AK13").Value = 10
$ws.Range("G16").Value = 3.9
$ws.Range("H16").Value = 4.2
$ws.Range("I16").Value = 1.62
$ws.Range("J16").Value = 5
$ws.Range("L16").Value = 2.25
$ws.Range("M16").Value = 1.03
$ws.Range("N16").Value = 15
$ws.Range("Q16").Value = 1.7
$ws.Range("R16").Value = 2.1
$ws.Range("W16").Value = 13
$ws.Range("X16").Value = 23
$ws.Range("Z16").Value = 51
$ws.Range("AA16").Value = 34
$ws.Range("AC16").Value = 15
$ws.Range("AD16").Value = 8.5
$ws.Range("AI16").Value = 8.5
$ws.Range("AK16").Value = 13
$ws.Range("AN16").Value = 6.5
$ws.Range("AQ16").Value = 81
$ws.Range("AR16").Value = 101
$ws.Range("AS16").Value = 201
$ws.Range("AY16").Value = 8.5
$ws.Range("AZ16").Value = 17
$ws.Range("BA16").Value = 26
$ws.Range("G17").Value = 2.05
$ws.Range("H17").Value = 3.5
$ws.Range("I17").Value = 3.6
$ws.Range("X17").Value = 9.5
$ws.Range("AI17").Value = 19
$ws.Range("AO17").Value = 11
$ws.Range("AY17").Value = 21
$ws.Range("G18").Value = 1.95
$ws.Range("H18").Value = 3
$ws.Range("I18").Value = 3.7
$ws.Range("J18").Value = 2.88
$ws.Range("K18").Value = 1.83
$ws.Range("L18").Value = 5
$ws.Range("M18").Value = 1.11
$ws.Range("N18").Value = 6.5
$ws.Range("O18").Value = 1.62
$ws.Range("Q18").Value = 2.88
$ws.Range("S18").Value = 1.62
$ws.Range("T18").Value = 2.2
$ws.Range("U18").Value = 2.38
$ws.Range("V18").Value = 1.53
$ws.Range("W18").Value = 5
$ws.Range("X18").Value = 8
$ws.Range("Y18").Value = 10
$ws.Range("Z18").Value = 17
$ws.Range("AA18").Value = 21
$ws.Range("AB18").Value = 41
$ws.Range("AC18").Value = 6
$ws.Range("AD18").Value = 6.5
$ws.Range("AE18").Value = 23
$ws.Range("AF18").Value = 101
$ws.Range("AH18").Value = 7.5
$ws.Range("AI18").Value = 17
$ws.Range("AJ18").Value = 15
$ws.Range("AK18").Value = 41
$ws.Range("AL18").Value = 41
$ws.Range("AM18").Value = 51
$ws.Range("AN18").Value = 3.75
$ws.Range("AO18").Value = 12
$ws.Range("AP18").Value = 29
$ws.Range("AQ18").Value = 41
$ws.Range("AR18").Value = 81
$ws.Range("AS18").Value = 351
$ws.Range("AT18").Value = 2.2
$ws.Range("AU18").Value = 10
$ws.Range("AV18").Value = 81
$ws.Range("AY18").Value = 26
$ws.Range("AZ18").Value = 41
$ws.Range("BA18").Value = 101
$ws.Range("BB18").Value = 151
$ws.Range("G21").Value = 1.45
$ws.Range("H21").Value = 5.5
$ws.Range("N21").Value = 26
$ws.Range("Q21").Value = 1.29
$ws.Range("R21").Value = 3.6
$ws.Range("X21").Value = 12
$ws.Range("AB21").Value = 15
$ws.Range("AC21").Value = 34
$ws.Range("AI21").Value = 41
$ws.Range("AN21").Value = 4.33
$ws.Range("AP21").Value = 12
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 3.9
$ws.Range("I22").Value = 2.1
$ws.Range("K22").Value = 2.38
$ws.Range("N22").Value = 17
$ws.Range("Q22").Value = 1.6
$ws.Range("R22").Value = 2.3
$ws.Range("U22").Value = 1.53
$ws.Range("V22").Value = 2.38
$ws.Range("W22").Value = 13
$ws.Range("AB22").Value = 23
$ws.Range("AC22").Value = 17
$ws.Range("AD22").Value = 8
$ws.Range("AE22").Value = 13
$ws.Range("AH22").Value = 11
$ws.Range("AO22").Value = 17
$ws.Range("AW22").Value = 351
$ws.Range("AY22").Value = 11
$ws.Range("AZ22").Value = 17
$ws.Range("BA22").Value = 34
$ws.Range("G23").Value = 2.3
$ws.Range("H23").Value = 4
$ws.Range("I23").Value = 2.63
$ws.Range("K23").Value = 2.63
$ws.Range("N23").Value = 23
$ws.Range("O23").Value = 1.08
$ws.Range("P23").Value = 8
$ws.Range("Q23").Value = 1.33
$ws.Range("R23").Value = 3.4
$ws.Range("S23").Value = 1.2
$ws.Range("T23").Value = 4.33
$ws.Range("U23").Value = 1.3
$ws.Range("V23").Value = 3.4
$ws.Range("W23").Value = 19
$ws.Range("X23").Value = 19
$ws.Range("Y23").Value = 11
$ws.Range("AA23").Value = 15
$ws.Range("AC23").Value = 29
$ws.Range("AD23").Value = 9.5
$ws.Range("AE23").Value = 10
$ws.Range("AF23").Value = 23
$ws.Range("AG23").Value = 51
$ws.Range("AH23").Value = 19
$ws.Range("AI23").Value = 21
$ws.Range("AT23").Value = 4.33
$ws.Range("AW23").Value = 151
$ws.Range("AZ23").Value = 15
$ws.Range("BC23").Value = 67
$ws.Range("G25").Value = 2.15
$ws.Range("I25").Value = 3.1
$ws.Range("J25").Value = 3
$ws.Range("L25").Value = 4
$ws.Range("Q25").Value = 2.25
$ws.Range("R25").Value = 1.62
$ws.Range("U25").Value = 1.91
$ws.Range("V25").Value = 1.8
$ws.Range("X25").Value = 10
$ws.Range("Z25").Value = 21
$ws.Range("AA25").Value = 21
$ws.Range("AF25").Value = 51
$ws.Range("AG25").Value = 900
$ws.Range("AI25").Value = 15
$ws.Range("AJ25").Value = 12
$ws.Range("AK25").Value = 34
$ws.Range("AL25").Value = 29
$ws.Range("AO25").Value = 13
$ws.Range("AX25").Value = 5
$ws.Range("AY25").Value = 19
$ws.Range("AZ25").Value = 29
$ws.Range("BC25").Value = 251
$ws.Range("G26").Value = 1.95
$ws.Range("I26").Value = 3.5
$ws.Range("K26").Value = 2.1
$ws.Range("L26").Value = 4.33
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 10
$ws.Range("O26").Value = 1.33
$ws.Range("P26").Value = 3.25
$ws.Range("Q26").Value = 2.08
$ws.Range("R26").Value = 1.73
$ws.Range("W26").Value = 7
$ws.Range("X26").Value = 9
$ws.Range("AG26").Value = 351
$ws.Range("AN26").Value = 4
$ws.Range("AS26").Value = 151
$ws.Range("AV26").Value = 51
$ws.Range("AY26").Value = 21
$ws.Range("AZ26").Value = 29
$ws.Range("BA26").Value = 67
$ws.Range("G29").Value = 2.35
$ws.Range("N29").Value = 7.5
$ws.Range("AL29").Value = 26
$ws.Range("AX29").Value = 4.75
$ws.Range("AY29").Value = 17
$ws.Range("O30").Value = 1.3
$ws.Range("P30").Value = 3.4
$ws.Range("Q30").Value = 2.05
$ws.Range("R30").Value = 1.75
$ws.Range("G33").Value = 1.8
$ws.Range("H33").Value = 3.25
$ws.Range("I33").Value = 4.5
$ws.Range("L33").Value = 4.85
$ws.Range("N33").Value = 9.3
$ws.Range("Q33").Value = 2.02
$ws.Range("R33").Value = 1.62
$ws.Range("U33").Value = 1.87
$ws.Range("V33").Value = 1.75
$ws.Range("W33").Value = 6.2
$ws.Range("X33").Value = 8
$ws.Range("Z33").Value = 15
$ws.Range("AC33").Value = 8
$ws.Range("AD33").Value = 6.3
$ws.Range("AE33").Value = 16.5
$ws.Range("AH33").Value = 10.5
$ws.Range("AM33").Value = 60
$ws.Range("AP33").Value = 17.5
$ws.Range("AR33").Value = 60
$ws.Range("AS33").Value = 250
$ws.Range("AT33").Value = 2.5
$ws.Range("AU33").Value = 7.2
$ws.Range("AV33").Value = 65
$ws.Range("AX33").Value = 6.2
$ws.Range("AY33").Value = 27
$ws.Range("BA33").Value = 175
